# cs-en-us-063pct.xlsx weekly refresh: new crime data collected.
# Bumps the report header (volume number / week-covering dates) and
# rewrites the CompStat precinct table's weekly/28-day/YTD/2-year figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a "no data" marker cell. These display either "0" or
# "***.*" and are stored as shared-string text (not numbers) in this
# workbook, matching the style already used by the neighboring cells
# on row 14 (C14/D14/G14 hold "0", E14/H14 hold "***.*").
# ---------------------------------------------------------------------
function Set-StringCell($ws, $cellRef, $text) {
    $sourceRef = "C14"
    if ($text -eq "***.*") { $sourceRef = "E14" }
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    $ws.Range($sourceRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------
# Helper: turn a former "no data" text marker cell back into a real
# number, picking up the plain numeric style from a neighboring cell
# that already carries it.
# ---------------------------------------------------------------------
function Set-NumericCell($ws, $cellRef, $val) {
    $ws.Range("D28").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null    # xlPasteFormats
    $ws.Range($cellRef).Value = $val
}

# ---------------------------------------------------------------------
# Header: "Volume 32   Number  27" -> "...  28"
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "28"

# ---------------------------------------------------------------------
# Header: reporting week "6/30/2025 ... 7/6/2025" -> "7/7/2025 ... 7/13/2025"
# ---------------------------------------------------------------------
$ws.Range("C9").Characters(27, 9).Text = "7/7/2025"
$afterFirst = $ws.Range("C9").Value2
$secondIdx = $afterFirst.IndexOf("7/6/2025")
$ws.Range("C9").Characters($secondIdx + 1, 8).Text = "7/13/2025"

# ---- String-type target cells (style 13, shared string "0"/"***.*") ----
Set-StringCell $ws "F14" "0"
Set-StringCell $ws "C15" "0"
Set-StringCell $ws "D15" "0"
Set-StringCell $ws "E15" "***.*"
Set-StringCell $ws "C16" "0"
Set-StringCell $ws "D16" "0"
Set-StringCell $ws "E16" "***.*"
Set-StringCell $ws "D23" "0"
Set-StringCell $ws "E23" "***.*"
Set-StringCell $ws "C27" "0"
Set-StringCell $ws "D27" "0"
Set-StringCell $ws "E27" "***.*"
Set-StringCell $ws "F31" "0"
Set-StringCell $ws "F33" "0"

# ---- Reverse (string -> numeric) cells ----
Set-NumericCell $ws "C28" 2

# ---- Plain numeric value changes ----
$ws.Range("N14").Value = -90.909090909090
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 20
$ws.Range("L16").Value = -8.771929824561
$ws.Range("M16").Value = -59.055118110236
$ws.Range("N16").Value = -85.674931129476
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -28.571428571428
$ws.Range("I17").Value = 78
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = -6.024096385542
$ws.Range("L17").Value = -14.285714285714
$ws.Range("M17").Value = 21.875
$ws.Range("N17").Value = -54.385964912280
$ws.Range("D18").Value = 2
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -66.666666666666
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = -2.857142857142
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -75
$ws.Range("N18").Value = -94.840667678300
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -7.894736842105
$ws.Range("I19").Value = 266
$ws.Range("J19").Value = 314
$ws.Range("K19").Value = -15.286624203821
$ws.Range("L19").Value = -19.393939393939
$ws.Range("M19").Value = 9.016393442622
$ws.Range("N19").Value = -24.216524216524
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 40
$ws.Range("I20").Value = 75
$ws.Range("J20").Value = 81
$ws.Range("K20").Value = -7.407407407407
$ws.Range("L20").Value = 15.384615384615
$ws.Range("M20").Value = -11.764705882352
$ws.Range("N20").Value = -95.003331112591
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 7.142857142857
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -8.108108108108
$ws.Range("I21").Value = 515
$ws.Range("J21").Value = 588
$ws.Range("K21").Value = -12.414965986394
$ws.Range("L21").Value = -14.309484193011
$ws.Range("M21").Value = -22.439759036144
$ws.Range("N21").Value = -83.235677083333
$ws.Range("L23").Value = -23.076923076923
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -27.586206896551
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = -40.718562874251
$ws.Range("I24").Value = 679
$ws.Range("J24").Value = 868
$ws.Range("K24").Value = -21.774193548387
$ws.Range("L24").Value = 6.593406593406
$ws.Range("M24").Value = 42.647058823529
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -36.363636363636
$ws.Range("F25").Value = 65
$ws.Range("G25").Value = 129
$ws.Range("H25").Value = -49.612403100775
$ws.Range("I25").Value = 482
$ws.Range("J25").Value = 717
$ws.Range("K25").Value = -32.775453277545
$ws.Range("L25").Value = 9.794988610478
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 66.666666666666
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = -11.538461538461
$ws.Range("I26").Value = 158
$ws.Range("J26").Value = 174
$ws.Range("K26").Value = -9.195402298850
$ws.Range("L26").Value = 13.669064748201
$ws.Range("M26").Value = -14.594594594594
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 4
$ws.Range("I28").Value = 15
$ws.Range("J28").Value = 17
$ws.Range("K28").Value = -11.764705882352
$ws.Range("L28").Value = -11.764705882352
$ws.Range("N29").Value = -82.142857142857
$ws.Range("N30").Value = -80
